# Update countries & provincias Spain
# Applies the 20-Jun-2020 13:33 -> 14:50 refresh of the COVID country table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp footer (row 1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Junio de 2020 a las 14:50"

# --- Two country names that swapped rank / position in the source list ----
# (Dominica <-> Fiyi, rows 202/203)
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"

# (Islas Turcas y Caicos <-> Santa Sede, rows 208/209)
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("A209").Value = "Islas Turcas y Caicos"

# --- Updated numeric stats per country row ---------------------------------
# row 4: Estados Unidos
$ws.Range("B4").Value = 2298108
$ws.Range("C4").Value = 918
$ws.Range("E4").Value = 1220607
$ws.Range("G4").Value = 17
$ws.Range("H4").Value = 121424

# row 7: India
$ws.Range("B7").Value = 396874
$ws.Range("C7").Value = 1062
$ws.Range("D7").Value = 214868
$ws.Range("E7").Value = 169034

# row 14: Alemania
$ws.Range("B14").Value = 190703
$ws.Range("C14").Value = 43
$ws.Range("E14").Value = 7343

# row 19: Arabia Saudita
$ws.Range("B19").Value = 154233
$ws.Range("C19").Value = 3941
$ws.Range("D19").Value = 98917
$ws.Range("E19").Value = 54086
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = 1230

# row 23: Catar
$ws.Range("B23").Value = 86488
$ws.Range("C23").Value = 1026
$ws.Range("D23").Value = 66763
$ws.Range("E23").Value = 19631
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 94

# row 31: Paises Bajos
$ws.Range("B31").Value = 49502
$ws.Range("C31").Value = 76
$ws.Range("G31").Value = 8
$ws.Range("H31").Value = 6089

# row 36: Kuwait
$ws.Range("B36").Value = 39145
$ws.Range("C36").Value = 467
$ws.Range("D36").Value = 30726
$ws.Range("E36").Value = 8100
$ws.Range("G36").Value = 6
$ws.Range("H36").Value = 319

# row 76: Uzbekistan
$ws.Range("B76").Value = 6119
$ws.Range("C76").Value = 173
$ws.Range("E76").Value = 1827

# row 101: Croacia
$ws.Range("B101").Value = 2299
$ws.Range("C101").Value = 19
$ws.Range("E101").Value = 50

# row 108: Albania
$ws.Range("E108").Value = 722
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 43

# row 208: Islas Turcas y Caicos -> Santa Sede (value text already updated above)
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

# row 209: Santa Sede -> Islas Turcas y Caicos (value text already updated above)
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1
